$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTT Test Hour Log")

# Fill in row 10 with new hour-log entry (Description first, so shared-string
# index order matches: 19 = Description text, 20 = Subject text)
$ws.Range("D10").Value = "Realized that with the current cell collection checking where neighbours `nis really annoying. So I started work on a cell grid that is made from the`ncollection of cells. This grid will be used to check possible empty `nneighbours, making the system ready for pretty much every algorithm."
$ws.Range("A10").Value = "Creating a cell grid from collection of cells."
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "5/22/2024"

# D10 should use the wrapped-text description style like the other rows
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").WrapText = $true

# Match row height of similar wrapped-text rows
$ws.Rows.Item(10).RowHeight = 52.5

# Update selection to reflect the edited cell
$ws.Range("F10").Select()
